$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The edit appends a new "For Presentation" / gh-pages workflow section
# at the end of the document (after "Then ($ cd link to folder to make
# it new master branch)"), and moves the "_GoBack" bookmark (Word's
# last-edit-location marker) so it ends up wrapping the very end of the
# newly typed content, matching what Word does when new text is typed
# at the end of a document.
# ------------------------------------------------------------------

# Remove the old _GoBack bookmark first - it will be re-created at the
# new insertion point once the new text has been typed.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# New lines to append, in order. The last one receives the _GoBack
# bookmark (collapsed, right after its text) once it has been typed.
$newLines = @(
    "For Presentation",
    "`$ git branch gh-pages",
    "`$ git checkout gh-pages",
    "`$ ls",
    "`$ git push origin gh-pages",
    "`$ touch .nojekyll",
    "`$ git add .nojekyll",
    "`$ git commit –a -m “Added a .nojekyll file”",
    "`$ git push origin gh-pages"
)

for ($i = 0; $i -lt $newLines.Count; $i++) {
    # Always type at the very end of the document.
    $insPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
    $insPoint.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newRange = $newPara.Range
    # The freshly created paragraph mark can inherit character formatting
    # (e.g. the yellow highlight) from whatever preceded it - make sure
    # the new text is plain/unhighlighted, as in the source document.
    $newRange.HighlightColorIndex = 0
    $newRange.InsertAfter($newLines[$i])

    if ($i -eq $newLines.Count - 1) {
        # Re-home "_GoBack" at the end of the very last content line,
        # while that paragraph is still the last paragraph in the
        # document (required for the bookmark to land correctly).
        $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
        $bmRange = $lastPara.Range.Duplicate
        $bmRange.MoveEnd(1, -1)
        $bmRange.Collapse(0)
        $bmRange.HighlightColorIndex = 0
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}

# Trailing blank paragraph that closes out the document body.
$tail = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$tail.InsertParagraphAfter()

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
